$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header cells — new test-case variants of "google.com"
$ws.Range("L1").Value = " google.com"
$ws.Range("M1").Value = "google.com "
$ws.Range("N1").Value = "google .com"
$ws.Range("O1").Value = "google. Com"

# Row 2 data marks ("x") for the new columns
$ws.Range("L2").Value = "x"
$ws.Range("M2").Value = "x"
$ws.Range("N2").Value = "x"
$ws.Range("O2").Value = "x"

# Move the active selection to L2 (bottom-right frozen pane)
$null = $ws.Range("L2").Select()
